$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Guca2a"
$ws.Range("C2").Value = "Gucy2c"
$ws.Range("D2").Value = "Resolving-Mac"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.7249690000000001
$ws.Range("H2").Value = 2.174907
$ws.Range("I2").Value = 0.3184714299144
$ws.Range("J2").Value = 0.3184714299144
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.015763
$ws.Range("N2").Value = 0.047289
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.011427686347
$ws.Range("R2").Value = 0.102849177123
$ws.Range("S2").Value = 0.3184714299144
$ws.Range("T2").Value = 0.3184714299144

# Row 3
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Guca2a"
$ws.Range("C3").Value = "Gucy2c"
$ws.Range("D3").Value = "Resolving-Mac"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.042549666666667
$ws.Range("H3").Value = 3.127649
$ws.Range("I3").Value = 0.4579813524441934
$ws.Range("J3").Value = 0.4579813524441934
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.015763
$ws.Range("N3").Value = 0.047289
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.01643371039566667
$ws.Range("R3").Value = 0.147903393561
$ws.Range("S3").Value = 0.4579813524441934
$ws.Range("T3").Value = 0.4579813524441934

# Row 4
$ws.Range("A4").Value = "MuSCs"
$ws.Range("B4").Value = "Guca2a"
$ws.Range("C4").Value = "Gucy2c"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.4282536666666667
$ws.Range("H4").Value = 1.284761
$ws.Range("I4").Value = 0.1881274338480931
$ws.Range("J4").Value = 0.188127433848093
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.015763
$ws.Range("N4").Value = 0.047289
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.006750562547666667
$ws.Range("R4").Value = 0.060755062929
$ws.Range("S4").Value = 0.1881274338480931
$ws.Range("T4").Value = 0.188127433848093

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("B5").Value = "Guca2a"
$ws.Range("C5").Value = "Gucy2c"
$ws.Range("D5").Value = "Resolving-Mac"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.08062966666666667
$ws.Range("H5").Value = 0.241889
$ws.Range("I5").Value = 0.0354197837933136
$ws.Range("J5").Value = 0.0354197837933136
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.015763
$ws.Range("N5").Value = 0.047289
$ws.Range("O5").Value = 1
$ws.Range("P5").Value = 1
$ws.Range("Q5").Value = 0.001270965435666667
$ws.Range("R5").Value = 0.011438688921
$ws.Range("S5").Value = 0.0354197837933136
$ws.Range("T5").Value = 0.0354197837933136
